# Updated cryptos list with GitHub Actions - refresh price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.995.31"
$ws.Range("E2").Value = "'  +1.02%  "
$ws.Range("D3").Value = "'1.641.59"
$ws.Range("E3").Value = "'  +0.48%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'212.90"
$ws.Range("E5").Value = "'  +0.35%  "
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'23.52"
$ws.Range("E8").Value = "'  +1.34%  "
$ws.Range("E9").Value = "'  -1.97%  "
$ws.Range("E10").Value = "'  +0.46%  "
$ws.Range("E11").Value = "'  +2.49%  "
$ws.Range("D12").Value = "'1.874.34"
$ws.Range("E12").Value = "'  +0.50%  "
$ws.Range("D13").Value = "'1.655.83"
$ws.Range("E13").Value = "'  +1.40%  "
$ws.Range("E14").Value = "'  +3.76%  "
$ws.Range("E15").Value = "'  +1.38%  "
$ws.Range("D16").Value = "'65.88"
$ws.Range("E16").Value = "'  +1.05%  "
$ws.Range("D17").Value = "'27.997.59"
$ws.Range("E17").Value = "'  +1.18%  "
$ws.Range("D18").Value = "'236.41"
$ws.Range("E18").Value = "'  +2.81%  "
$ws.Range("D19").Value = "'0.0₃0725"
$ws.Range("E19").Value = "'  +0.63%  "
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "'  +0.73%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "'  +0.16%  "
$ws.Range("E23").Value = "'  +0.79%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "'  -2.05%  "
$ws.Range("D25").Value = "'151.78"
$ws.Range("E25").Value = "'  +1.92%  "
$ws.Range("D26").Value = "'6.96"
$ws.Range("E26").Value = "'  +1.37%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("E27").Value = "'  +0.70%  "
$ws.Range("E28").Value = "'  +0.13%  "
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("E30").Value = "'  +0.46%  "
$ws.Range("E31").Value = "'  +0.61%  "
$ws.Range("E32").Value = "'  +1.86%  "
$ws.Range("D33").Value = "'3.12"
$ws.Range("E33").Value = "'  +1.40%  "
$ws.Range("D34").Value = "'1.417.37"
$ws.Range("E34").Value = "'  -3.76%  "
$ws.Range("E35").Value = "'  +2.58%  "
$ws.Range("E36").Value = "'  +1.45%  "
$ws.Range("E37").Value = "'  +1.74%  "
$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "'  +0.72%  "
$ws.Range("D39").Value = "'0.558"
$ws.Range("E39").Value = "'  -0.17%  "
$ws.Range("D40").Value = "'0.899"
$ws.Range("E40").Value = "'  -4.74%  "
$ws.Range("E41").Value = "'  +1.16%  "
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("D43").Value = "'1.86"
$ws.Range("E43").Value = "'  +6.57%  "
$ws.Range("D44").Value = "'66.62"
$ws.Range("E44").Value = "'  -1.78%  "
$ws.Range("E45").Value = "'  +3.05%  "
$ws.Range("E46").Value = "'  -0.05%  "
$ws.Range("D47").Value = "'1.783.44"
$ws.Range("E47").Value = "'  +0.63%  "
$ws.Range("D48").Value = "'87.91"
$ws.Range("E49").Value = "'  +0.80%  "
$ws.Range("E50").Value = "'  +0.32%  "
$ws.Range("D51").Value = "'7.61"
$ws.Range("E51").Value = "'  -1.25%  "
